# Scheduled market-price refresh for the Leve profit-tracker sheets.
# Updates currentAveragePrice / NQ / HQ columns (H-N) for the leves whose
# market snapshot changed, per class-sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1533.0834
$ws.Cells.Item(17, 10).Value = 1533.0834
$ws.Cells.Item(17, 12).Value = 4599.2502
$ws.Cells.Item(17, 14).Value = -4935.2502
$ws.Cells.Item(32, 8).Value = 929.6539
$ws.Cells.Item(32, 9).Value = 798.7368
$ws.Cells.Item(32, 11).Value = 798.7368
$ws.Cells.Item(32, 13).Value = -472.7368
$ws.Cells.Item(132, 8).Value = 1841.7142
$ws.Cells.Item(132, 9).Value = 1832
$ws.Cells.Item(132, 11).Value = 5496
$ws.Cells.Item(132, 13).Value = -2966
$ws.Cells.Item(137, 8).Value = 4978
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 4978
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).Value = -20034
$ws.Cells.Item(138, 8).Value = 2777.7778
$ws.Cells.Item(138, 9).Value = 1000
$ws.Cells.Item(138, 10).Value = 3000
$ws.Cells.Item(138, 11).Value = 3000
$ws.Cells.Item(138, 12).Value = 9000
$ws.Cells.Item(138, 13).Value = 2140
$ws.Cells.Item(138, 14).Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 12143.9
$ws.Cells.Item(2, 9).Value = 1070
$ws.Cells.Item(2, 10).Value = 37983
$ws.Cells.Item(2, 11).Value = 1070
$ws.Cells.Item(2, 12).Value = 37983
$ws.Cells.Item(2, 13).Value = -957
$ws.Cells.Item(2, 14).Value = -38209
$ws.Cells.Item(5, 8).Value = 122.9
$ws.Cells.Item(5, 9).Value = 111.625
$ws.Cells.Item(5, 10).Value = 168
$ws.Cells.Item(5, 11).Value = 111.625
$ws.Cells.Item(5, 12).Value = 168
$ws.Cells.Item(5, 13).Value = 0.375
$ws.Cells.Item(5, 14).Value = -392
$ws.Cells.Item(61, 8).Value = 3236.5
$ws.Cells.Item(61, 9).Value = 2982.5
$ws.Cells.Item(61, 11).Value = 2982.5
$ws.Cells.Item(61, 13).Value = -2770.5
$ws.Cells.Item(74, 8).Value = 1916.5
$ws.Cells.Item(74, 9).Value = 1899.8
$ws.Cells.Item(74, 10).Value = 2000
$ws.Cells.Item(74, 11).Value = 1899.8
$ws.Cells.Item(74, 12).Value = 2000
$ws.Cells.Item(74, 13).Value = -1025.8
$ws.Cells.Item(74, 14).Value = -3748
$ws.Cells.Item(77, 8).Value = 1916.5
$ws.Cells.Item(77, 9).Value = 1899.8
$ws.Cells.Item(77, 10).Value = 2000
$ws.Cells.Item(77, 11).Value = 9499
$ws.Cells.Item(77, 12).Value = 10000
$ws.Cells.Item(77, 13).Value = -5131
$ws.Cells.Item(77, 14).Value = -18736
$ws.Cells.Item(102, 8).Value = 2029.4445
$ws.Cells.Item(102, 9).Value = 1630
$ws.Cells.Item(102, 11).Value = 1630
$ws.Cells.Item(102, 13).Value = -8
$ws.Cells.Item(116, 8).Value = 12143.9
$ws.Cells.Item(116, 9).Value = 1070
$ws.Cells.Item(116, 10).Value = 37983
$ws.Cells.Item(116, 11).Value = 1070
$ws.Cells.Item(116, 12).Value = 37983
$ws.Cells.Item(116, 13).Value = 1224
$ws.Cells.Item(116, 14).Value = -42571
$ws.Cells.Item(135, 8).Value = 29996.334
$ws.Cells.Item(135, 10).Value = 29996.334
$ws.Cells.Item(135, 12).Value = 29996.334
$ws.Cells.Item(135, 14).Value = -40136.334
$ws.Cells.Item(136, 8).Value = 3236.5
$ws.Cells.Item(136, 9).Value = 2982.5
$ws.Cells.Item(136, 11).Value = 8947.5
$ws.Cells.Item(136, 13).Value = -6397.5
$ws.Cells.Item(139, 8).Value = 49998.5
$ws.Cells.Item(139, 10).Value = 49998.5
$ws.Cells.Item(139, 12).Value = 49998.5
$ws.Cells.Item(139, 14).Value = -60278.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 12143.9
$ws.Cells.Item(3, 9).Value = 1070
$ws.Cells.Item(3, 10).Value = 37983
$ws.Cells.Item(3, 11).Value = 1070
$ws.Cells.Item(3, 12).Value = 37983
$ws.Cells.Item(3, 13).Value = -956
$ws.Cells.Item(3, 14).Value = -38211
$ws.Cells.Item(4, 8).Value = 122.9
$ws.Cells.Item(4, 9).Value = 111.625
$ws.Cells.Item(4, 10).Value = 168
$ws.Cells.Item(4, 11).Value = 111.625
$ws.Cells.Item(4, 12).Value = 168
$ws.Cells.Item(4, 13).Value = 3.375
$ws.Cells.Item(4, 14).Value = -398
$ws.Cells.Item(15, 8).Value = 29976.223
$ws.Cells.Item(15, 9).Value = 18000
$ws.Cells.Item(15, 10).Value = 35964.332
$ws.Cells.Item(15, 11).Value = 18000
$ws.Cells.Item(15, 12).Value = 35964.332
$ws.Cells.Item(15, 13).Value = -17773
$ws.Cells.Item(15, 14).Value = -36418.332
$ws.Cells.Item(19, 8).Value = 25979
$ws.Cells.Item(19, 10).Value = 25979
$ws.Cells.Item(19, 12).Value = 25979
$ws.Cells.Item(19, 14).Value = -26325
$ws.Cells.Item(86, 8).Value = 4501.2
$ws.Cells.Item(86, 9).Value = 4876.5
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 4876.5
$ws.Cells.Item(86, 12).Value = 3000
$ws.Cells.Item(86, 13).Value = -3753.5
$ws.Cells.Item(86, 14).Value = -5246
$ws.Cells.Item(89, 8).Value = 4501.2
$ws.Cells.Item(89, 9).Value = 4876.5
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 24382.5
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = -18766.5
$ws.Cells.Item(89, 14).Value = -26232
$ws.Cells.Item(94, 8).Value = 2442.111
$ws.Cells.Item(94, 9).Value = 2442.111
$ws.Cells.Item(94, 11).Value = 2442.111
$ws.Cells.Item(94, 13).Value = -1991.111
$ws.Cells.Item(134, 8).Value = 6656.8887
$ws.Cells.Item(134, 9).Value = 6656.8887
$ws.Cells.Item(134, 11).Value = 19970.6661
$ws.Cells.Item(134, 13).Value = -17435.6661
$ws.Cells.Item(138, 8).Value = 43332.332
$ws.Cells.Item(138, 10).Value = 43332.332
$ws.Cells.Item(138, 12).Value = 43332.332
$ws.Cells.Item(138, 14).Value = -53612.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 3796
$ws.Cells.Item(132, 9).Value = 3745
$ws.Cells.Item(132, 11).Value = 11235
$ws.Cells.Item(132, 13).Value = -8705

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 364.9
$ws.Cells.Item(47, 9).Value = 321.2857
$ws.Cells.Item(47, 10).Value = 466.66666
$ws.Cells.Item(47, 11).Value = 963.8571000000001
$ws.Cells.Item(47, 12).Value = 1399.99998
$ws.Cells.Item(47, 13).Value = -532.8571000000001
$ws.Cells.Item(47, 14).Value = -2261.99998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1750
$ws.Cells.Item(132, 9).Value = 1000
$ws.Cells.Item(132, 10).Value = 2500
$ws.Cells.Item(132, 11).Value = 3000
$ws.Cells.Item(132, 12).Value = 7500
$ws.Cells.Item(132, 13).Value = -470
$ws.Cells.Item(132, 14).Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 957
$ws.Cells.Item(22, 9).Value = 857.6
$ws.Cells.Item(22, 11).Value = 857.6
$ws.Cells.Item(22, 13).Value = -562.6
$ws.Cells.Item(27, 8).Value = 957
$ws.Cells.Item(27, 9).Value = 857.6
$ws.Cells.Item(27, 11).Value = 857.6
$ws.Cells.Item(27, 13).Value = -750.6
$ws.Cells.Item(40, 8).Value = 6487.5
$ws.Cells.Item(40, 9).Value = 6342.857
$ws.Cells.Item(40, 10).Value = 7500
$ws.Cells.Item(40, 11).Value = 6342.857
$ws.Cells.Item(40, 12).Value = 7500
$ws.Cells.Item(40, 13).Value = -6206.857
$ws.Cells.Item(40, 14).Value = -7772
